$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sources worksheet: reorder the first six header columns, add the new
# sample data rows 2-12, and refresh the autofilter / defined name.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sources")

# Reorder header row A1:F1 (columns G:L are unchanged)
$ws.Range("A1").Value = "Statistics_Source_Name"
$ws.Range("B1").Value = "Statistics_Source_Retrieval_Code"
$ws.Range("C1").Value = "``Stats`` Vendor_ID"
$ws.Range("D1").Value = "Statistics_Source_ID"
$ws.Range("E1").Value = "Resource_Source_ID"
$ws.Range("F1").Value = "Resource_Source_Name"

# New sample data rows (only columns A, C, D are populated)
$names = @(
    "ProQuest",
    "EBSCOhost",
    "Gale Cengage Learning",
    "iG Library/Business Expert Press (BEP)",
    "DemographicsNow",
    "Ebook Central",
    "Peterson's Career Prep",
    "Peterson's Test Prep",
    "Peterson's Prep",
    "Pivot",
    "UlrichsWeb"
)
$statsVendorIds = @(1, 2, 3, 4, 3, 1, 3, 3, 3, 1, 1)
$statsSourceIds = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $statsVendorIds[$i]
    $ws.Cells.Item($row, 4).Value = $statsSourceIds[$i]
}

# Recompute best-fit column widths now that the content has changed
$ws.Columns.AutoFit()

# Refresh the autofilter so it covers the new data extent
$ws.AutoFilterMode = $false
$ws.Range("A1:L12").AutoFilter(1)

# Keep the "_FilterDatabase" defined name for Sources in sync
foreach ($dn in $wb.Names) {
    if ($dn.Name -eq "Sources!_FilterDatabase") {
        $dn.RefersTo = "=Sources!`$A`$1:`$L`$12"
    }
}

# ---------------------------------------------------------------------------
# initialize_vendors worksheet: restore its own (now unselected) selection
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("initialize_vendors")
$ws2.Activate()
$ws2.Range("D7").Select()

# Update the selection/active cell and make Sources the selected tab (last,
# so it ends up as the workbook's active sheet)
$ws.Activate()
$ws.Range("F3").Select()
